# Regression_Models_CD_Training_Metrics.xlsx edit
# - Reorders the model rows so each row carries its own "Training Set (320 Compounds)"
#   label in column A (the old A3:A9 merge is removed) and every model gets its own row
#   (Linear Regression, Support Vector Regression, K-Nearest Neighbour Regressor,
#    Random Forest Regressor, Decision Tree Regressor, Stochastic Gradient Descent
#    Regressor all shift up a row because the old blank-merged rows 4-9 are now filled).
# - Updates formatting: numeric format "0.00000", header row no longer vertically+horizontally
#   centered for the "Set" column, row heights, and un-merges A3:A9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Unmerge the old A3:A9 "Training Set (320 Compounds)" block so each row
#    gets its own value in column A.
# ---------------------------------------------------------------------------
$ws.Range("A3:A9").UnMerge()

# ---------------------------------------------------------------------------
# 2. Rewrite the data rows (3-9) with the new model ordering / values.
# ---------------------------------------------------------------------------
$trainingSet = "Training Set (320 Compounds)"
$pValue = [double]"9.9009900990098994E-3"

$ws.Range("A3").Value = $trainingSet
$ws.Range("B3").Value = "Dummy Regressor"
$ws.Range("C3").Value = -0.62085023504635695
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "NA"

$ws.Range("A4").Value = $trainingSet
$ws.Range("B4").Value = "Linear Regression"
$ws.Range("C4").Value = -0.50098642309722496
$ws.Range("D4").Value = 0.19939900753380299
$ws.Range("E4").Value = $pValue

$ws.Range("A5").Value = $trainingSet
$ws.Range("B5").Value = "Support Vector Regression"
$ws.Range("C5").Value = -0.46663119180065998
$ws.Range("D5").Value = 0.26939204165821501
$ws.Range("E5").Value = $pValue

$ws.Range("A6").Value = $trainingSet
$ws.Range("B6").Value = "K-Nearest Neighbour Regressor"
$ws.Range("C6").Value = -0.50633411595638
$ws.Range("D6").Value = 0.227194343699572
$ws.Range("E6").Value = $pValue

$ws.Range("A7").Value = $trainingSet
$ws.Range("B7").Value = "Random Forest Regressor"
$ws.Range("C7").Value = -0.49887430664758098
$ws.Range("D7").Value = 0.12504893944994599
$ws.Range("E7").Value = $pValue

$ws.Range("A8").Value = $trainingSet
$ws.Range("B8").Value = "Decision Tree Regressor"
$ws.Range("C8").Value = -0.62121875299999996
$ws.Range("D8").Value = -1.25635869198143
$ws.Range("E8").Value = 0.16831699999999999

$ws.Range("A9").Value = $trainingSet
$ws.Range("B9").Value = "Stochastic Gradient Descent Regressor"
$ws.Range("C9").Value = -0.49768864453836198
$ws.Range("D9").Value = 0.26195709902844699
$ws.Range("E9").Value = $pValue

# ---------------------------------------------------------------------------
# 3. Formatting.
# ---------------------------------------------------------------------------

# Numeric format for all number cells becomes "0.00000" (was "0.0000000000").
$ws.Range("C3:D9").NumberFormat = "0.00000"
$ws.Range("E3:E9").NumberFormat = "0.00000"

# Header row ("Set","Model","Negated Mean Absolute Error ","R2 Score","Permutation
# Testing P-Value"): "Set"/"Model"/"R2 Score" (A2,B2,D2) already have the correct
# (no-alignment) formatting and are left untouched. "Negated Mean Absolute Error "
# and "Permutation Testing P-Value" (C2,E2) lose the horizontal=center/vertical=center
# alignment they used to have, keeping only wrapText.
$ws.Range("C2").HorizontalAlignment = 1   # xlGeneral
$ws.Range("C2").VerticalAlignment = 1     # xlGeneral
$ws.Range("E2").HorizontalAlignment = 1
$ws.Range("E2").VerticalAlignment = 1

# Column A (Training Set) keeps its vertical=center + wrapText formatting; only
# the horizontal=center alignment is dropped (back to general).
$ws.Range("A3:A9").HorizontalAlignment = 1  # xlGeneral

# Column B (model names) already uses left/general + no wrap - no change needed.

# Numeric columns C & D, and the "NA"/0.168317 cells in column E (E3 & E8) already
# have center alignment + wrap text - only the number format changed (handled above).

# The remaining E column p-value cells (E4,E5,E6,E7,E9) stay centered but lose wrapText.
$ws.Range("E4").WrapText = $false
$ws.Range("E5").WrapText = $false
$ws.Range("E6").WrapText = $false
$ws.Range("E7").WrapText = $false
$ws.Range("E9").WrapText = $false

# ---------------------------------------------------------------------------
# 4. Row heights: row 3 stays 15.6 (now an explicit height), rows 4-9 grow to
#    31.2 because column A/B now wrap onto two lines per row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 15.6
$ws.Rows.Item(4).RowHeight = 31.2
$ws.Rows.Item(5).RowHeight = 31.2
$ws.Rows.Item(6).RowHeight = 31.2
$ws.Rows.Item(7).RowHeight = 31.2
$ws.Rows.Item(8).RowHeight = 31.2
$ws.Rows.Item(9).RowHeight = 31.2

# ---------------------------------------------------------------------------
# 5. Selection moves from G8 to E4.
# ---------------------------------------------------------------------------
$ws.Range("E4").Select()

Write-Host "edit applied"
